# Apply updated crypto price/volume snapshot values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as literal text
# (mirrors typing a leading apostrophe in Excel), preserving the sheet's
# existing convention of keeping Price/Volume columns as text cells.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
}

$ws.Range("D2").Value = '63.031.05'
$ws.Range("E2").Value = '  -2.43%  '

$ws.Range("D3").Value = '3.120.84'
$ws.Range("E3").Value = '  -1.25%  '

Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  -0.18%  '

Set-TextValue $ws.Range("D5") '591.83'
$ws.Range("E5").Value = '  -3.36%  '

Set-TextValue $ws.Range("D6") '136.50'
$ws.Range("E6").Value = '  -5.41%  '

$ws.Range("E7").Value = '  -0.16%  '

$ws.Range("D8").Value = '3.121.50'
$ws.Range("E8").Value = '  -1.07%  '

Set-TextValue $ws.Range("D9") '0.518'
$ws.Range("E9").Value = '  -1.40%  '

Set-TextValue $ws.Range("D10") '0.146'
$ws.Range("E10").Value = '  -3.15%  '

Set-TextValue $ws.Range("D11") '5.30'
$ws.Range("E11").Value = '  -1.20%  '

Set-TextValue $ws.Range("D12") '0.457'
$ws.Range("E12").Value = '  -3.26%  '

Set-TextValue $ws.Range("D13") '0.0000247'
$ws.Range("E13").Value = '  -3.37%  '

Set-TextValue $ws.Range("D14") '33.96'
$ws.Range("E14").Value = '  -4.72%  '

$ws.Range("D15").Value = '3.624.03'
$ws.Range("E15").Value = '  -1.51%  '

$ws.Range("E16").Value = '  +1.52%  '

$ws.Range("D17").Value = '63.043.74'
$ws.Range("E17").Value = '  -2.38%  '

$ws.Range("D18").Value = '3.113.73'
$ws.Range("E18").Value = '  -1.47%  '

Set-TextValue $ws.Range("D19") '6.68'
$ws.Range("E19").Value = '  -2.73%  '

Set-TextValue $ws.Range("D20") '474.72'
$ws.Range("E20").Value = '  -0.48%  '

Set-TextValue $ws.Range("D21") '14.17'
$ws.Range("E21").Value = '  -3.47%  '

Set-TextValue $ws.Range("D22") '0.694'
$ws.Range("E22").Value = '  -4.34%  '

Set-TextValue $ws.Range("D23") '7.62'
$ws.Range("E23").Value = '  -3.06%  '

Set-TextValue $ws.Range("D24") '87.19'
$ws.Range("E24").Value = '  +3.00%  '

Set-TextValue $ws.Range("D25") '13.01'
$ws.Range("E25").Value = '  -5.53%  '

$ws.Range("E26").Value = '  +0.08%  '

Set-TextValue $ws.Range("D27") '2.70'
$ws.Range("E27").Value = '  -3.86%  '

$ws.Range("B28").Value = 'NEARProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D28") '7.09'
$ws.Range("E28").Value = '  -4.74%  '

$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D29") '7.98'
$ws.Range("E29").Value = '  -6.91%  '

Set-TextValue $ws.Range("D30") '2.04'
$ws.Range("E30").Value = '  -3.16%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D31") '27.21'
$ws.Range("E31").Value = '  +2.67%  '

$ws.Range("B32").Value = 'FirstDigitalUSD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D32") '1.00'
$ws.Range("E32").Value = '  -0.13%  '

Set-TextValue $ws.Range("D33") '0.108'
$ws.Range("E33").Value = '  -11.22%  '

Set-TextValue $ws.Range("D34") '2.53'
$ws.Range("E34").Value = '  -4.82%  '

$ws.Range("E35").Value = '  -3.20%  '

Set-TextValue $ws.Range("D36") '5.84'
$ws.Range("E36").Value = '  -1.95%  '

Set-TextValue $ws.Range("D37") '51.95'
$ws.Range("E37").Value = '  -1.62%  '

$ws.Range("E38").Value = '  -4.80%  '

Set-TextValue $ws.Range("D39") '0.0386'
$ws.Range("E39").Value = '  -2.81%  '

Set-TextValue $ws.Range("D40") '419.17'
$ws.Range("E40").Value = '  -7.68%  '

$ws.Range("E41").Value = '  -1.11%  '

$ws.Range("B42").Value = 'Cosmos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D42") '8.23'
$ws.Range("E42").Value = '  -1.14%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D43") '2.71'
$ws.Range("E43").Value = '  -13.39%  '

$ws.Range("D44").Value = '2.861.05'
$ws.Range("E44").Value = '  +0.60%  '

Set-TextValue $ws.Range("D45") '0.256'
$ws.Range("E45").Value = '  -4.42%  '

Set-TextValue $ws.Range("D47") '2.10'
$ws.Range("E47").Value = '  -7.27%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D48") '25.42'
$ws.Range("E48").Value = '  -4.26%  '

$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue $ws.Range("D49") '2.29'
$ws.Range("E49").Value = '  -7.40%  '

$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D50") '0.113'
$ws.Range("E50").Value = '  -0.72%  '

Set-TextValue $ws.Range("D51") '119.16'
$ws.Range("E51").Value = '  -1.03%  '
